$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (CSC103): update duration(s) and end timestamp
$ws.Range("B3").Value = 41
$ws.Range("D3").Value = 45689.6969383912

# Row 4 (MAT111): update duration(s) and end timestamp
$ws.Range("B4").Value = 15
$ws.Range("D4").Value = 45689.69669267517
